$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 407 (existing rows 407-426 shift down to 410-429,
# preserving all their values/formatting).
$ws.Rows("407:409").Insert()

# Populate the 3 newly inserted rows with the new weekly price-report entries
# (Comercializadora del Agro de Limari - Pepino dulce, fecha 2022-06-02).

# Row 407: Calidad "Especial"
$ws.Cells.Item(407, 1).Value = 2
$ws.Cells.Item(407, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(407, 3).Value = "Coquimbo"
$ws.Cells.Item(407, 4).Value = 44714
$ws.Cells.Item(407, 5).Value = 4
$ws.Cells.Item(407, 6).Value = 100112043
$ws.Cells.Item(407, 7).Value = "Pepino dulce"
$ws.Cells.Item(407, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(407, 9).Value = "Especial"
$ws.Cells.Item(407, 10).Value = 400
$ws.Cells.Item(407, 11).Value = 12000
$ws.Cells.Item(407, 12).Value = 13000
$ws.Cells.Item(407, 13).Value = 12500
$ws.Cells.Item(407, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(407, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(407, 16).Value = 694
$ws.Cells.Item(407, 17).Value = 18
$ws.Cells.Item(407, 18).Value = "Hortaliza"

# Row 408: Calidad "Primera"
$ws.Cells.Item(408, 1).Value = 2
$ws.Cells.Item(408, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(408, 3).Value = "Coquimbo"
$ws.Cells.Item(408, 4).Value = 44714
$ws.Cells.Item(408, 5).Value = 4
$ws.Cells.Item(408, 6).Value = 100112043
$ws.Cells.Item(408, 7).Value = "Pepino dulce"
$ws.Cells.Item(408, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(408, 9).Value = "Primera"
$ws.Cells.Item(408, 10).Value = 440
$ws.Cells.Item(408, 11).Value = 10000
$ws.Cells.Item(408, 12).Value = 11000
$ws.Cells.Item(408, 13).Value = 10500
$ws.Cells.Item(408, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(408, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(408, 16).Value = 583
$ws.Cells.Item(408, 17).Value = 18
$ws.Cells.Item(408, 18).Value = "Hortaliza"

# Row 409: Calidad "Segunda"
$ws.Cells.Item(409, 1).Value = 2
$ws.Cells.Item(409, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(409, 3).Value = "Coquimbo"
$ws.Cells.Item(409, 4).Value = 44714
$ws.Cells.Item(409, 5).Value = 4
$ws.Cells.Item(409, 6).Value = 100112043
$ws.Cells.Item(409, 7).Value = "Pepino dulce"
$ws.Cells.Item(409, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(409, 9).Value = "Segunda"
$ws.Cells.Item(409, 10).Value = 500
$ws.Cells.Item(409, 11).Value = 7000
$ws.Cells.Item(409, 12).Value = 8000
$ws.Cells.Item(409, 13).Value = 7500
$ws.Cells.Item(409, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(409, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(409, 16).Value = 417
$ws.Cells.Item(409, 17).Value = 18
$ws.Cells.Item(409, 18).Value = "Hortaliza"
